# =====================================================================
# feat: add 2022-Q4 data
#
# 1. Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. right
#    before the existing "2022-Q3" sheet), populated with the Q4 fund
#    holdings table.
# 2. Insert a matching summary row at the top of the "总计" data table.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: create + position the new "2022-Q4" worksheet
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"
$newSheet.Move($q3Sheet)

# NB: worksheet proxy objects are resolved positionally, so any handle
# obtained before a sheet reorder can end up pointing at the wrong sheet
# afterwards. Re-resolve by name once the order is settled.
$ws = $wb.Worksheets.Item("2022-Q4")

# ---------------------------------------------------------------------
# Part 2: header row (bold, bordered, centered - matches the other
# quarterly sheets already in the workbook)
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$headerRange = $ws.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ---------------------------------------------------------------------
# Part 3: data rows. Columns B:G ("基金代码".."持有市值(亿元)") are stored
# as text in the source data (e.g. "000362", "26.05"), not numbers, so
# force the number format to Text before assigning - otherwise the host
# auto-coerces numeric-looking strings into real numbers (and mangles
# fund codes by dropping their leading zeros).
# ---------------------------------------------------------------------
$q4Data = @(
    @("000362", "国泰聚信价值优势灵活配置混合A", "26.05", "91.85", "4.66", "1.2139", 6),
    @("000363", "国泰聚信价值优势灵活配置混合C", "12.90", "91.85", "4.66", "0.6011", 6),
    @("020010", "国泰金牛创新混合", "13.29", "89.45", "3.43", "0.4558", 8),
    @("012173", "国泰兴泽优选一年持有期混合A", "8.18", "92.65", "5.07", "0.4147", 4),
    @("200006", "长城消费增值混合", "5.59", "92.96", "5.88", "0.3287", 1),
    @("012174", "国泰兴泽优选一年持有期混合C", "5.86", "92.65", "5.07", "0.2971", 4),
    @("013890", "国泰睿毅三年持有期混合A", "4.82", "90.04", "5.54", "0.2670", 1),
    @("007835", "国泰鑫睿混合", "7.62", "79.19", "3.35", "0.2553", 5),
    @("011466", "兴业医疗保健混合A", "4.60", "88.16", "4.92", "0.2263", 8),
    @("003516", "国泰融安多策略灵活配置混合A", "7.77", "78.85", "2.55", "0.1981", 3),
    @("005244", "国泰聚优价值灵活配置混合A", "4.39", "90.86", "3.10", "0.1361", 10),
    @("011467", "兴业医疗保健混合C", "2.12", "88.16", "4.92", "0.1043", 8),
    @("090020", "大成健康产业混合A", "2.49", "92.90", "3.82", "0.0951", 10),
    @("008619", "永赢医药健康股票C", "1.78", "89.12", "5.11", "0.0910", 8),
    @("012045", "大成医药健康股票A", "2.04", "92.92", "3.78", "0.0771", 10),
    @("005245", "国泰聚优价值灵活配置混合C", "2.10", "90.86", "3.10", "0.0651", 10),
    @("020023", "国泰事件驱动策略混合A", "2.19", "82.03", "2.92", "0.0639", 6),
    @("012880", "国泰景气优选混合A", "3.24", "75.90", "1.93", "0.0625", 9),
    @("008618", "永赢医药健康股票A", "0.58", "89.12", "5.11", "0.0296", 8),
    @("013891", "国泰睿毅三年持有期混合C", "0.45", "90.04", "5.54", "0.0249", 1),
    @("620002", "金元顺安成长动力混合", "0.35", "72.02", "3.54", "0.0124", 4),
    @("012046", "大成医药健康股票C", "0.27", "92.92", "3.78", "0.0102", 10),
    @("012881", "国泰景气优选混合C", "0.22", "75.90", "1.93", "0.0042", 9),
    @("014960", "国泰融安多策略灵活配置混合C", "0.14", "78.85", "2.55", "0.0036", 3),
    @("015921", "申万菱信国证2000指数增强A", "0.21", "94.00", "0.54", "0.0011", 1),
    @("015922", "申万菱信国证2000指数增强C", "0.08", "94.00", "0.54", "0.0004", 1),
    @("016060", "大成健康产业混合C", "0.01", "92.90", "3.82", "0.0004", 10),
    @("015592", "国泰事件驱动策略混合C", "0.01", "82.03", "2.92", "0.0003", 6),
)

$textCols = $ws.Range("B2:G" + ($q4Data.Length + 1))
$textCols.NumberFormat = "@"

for ($r = 0; $r -lt $q4Data.Length; $r++) {
    $row = $q4Data[$r]
    $excelRow = $r + 2
    $ws.Cells.Item($excelRow, 1).Value = $r
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($excelRow, $c + 2).Value = $row[$c]
    }
    $ws.Cells.Item($excelRow, 8).Value = $row[6]
}

Write-Output "2022-Q4 sheet populated with $($q4Data.Length) rows"

# ---------------------------------------------------------------------
# Part 4: update the "总计" (summary) sheet - insert a new row for
# 2022-Q4 above the existing 2022-Q3 row, shifting everything else
# down, and keep the running index in column A sequential.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$totalRows = @(
    @("2022-Q4", 28, 5.04),
    @("2022-Q3", 18, 3.8),
    @("2022-Q2", 6, 0.98),
    @("2022-Q1", 5, 0.62),
    @("2021-Q2", 2, 0.03)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $row = $totalRows[$r]
    $excelRow = $r + 2
    $total.Cells.Item($excelRow, 1).Value = $r
    $total.Cells.Item($excelRow, 2).Value = $row[0]
    $total.Cells.Item($excelRow, 3).Value = $row[1]
    $total.Cells.Item($excelRow, 4).Value = $row[2]
}

Write-Output "总计 sheet updated with $($totalRows.Length) rows"
